$d = $word.ActiveDocument

# Simple whole-paragraph text replacements. We match on the full paragraph
# text (stripped of the trailing paragraph-mark / cell-mark characters) so
# that "asdasda" and "asdasdas" (one a prefix of the other) cannot cross-match.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Name: kathiravan") {
        $p.Range.Text = "Name: sfdsf"
    }
    elseif ($t -eq "Age: asdasd") {
        $p.Range.Text = "Age: sdfsdfjlkj"
    }
    elseif ($t -eq "asdasda") {
        $p.Range.Text = "jsdlkfjlkdsj"
    }
    elseif ($t -eq "asdasdas") {
        $p.Range.Text = "kjlkjdlkfjl"
    }
}

# The numbered-list paragraph mixes text runs with manual line breaks
# (<w:br/> => ^l). Drop the third line entirely while renaming the first two.
$d.Content.Find.Execute("1. asdasd^l2. asdsad^l3. adasdsa", $false, $false, $false, $false, $false, $true, 1, $false, "1. kjsdlkfjk^l2. sdfdsf", 2)
